# "Generate Report for Handoff"
#
# The localization-status report moves from "In Translation" to
# "Ready for handoff": the status text + the handoff timestamps are
# refreshed, and the (auto-fit) Status/Handoff-datetime columns widen
# to fit the new, longer "Ready for handoff" label.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value     = "Ready for handoff"   # Status column
$wsDeDe.Range("C2").Value     = "Ready for handoff"   # Status column

# --- Handoff timestamps ------------------------------------------------
# (re-assert the datetime display format on the cells we touch so the
# existing "yyyy-mm-dd HH:mm:ss" formatting survives the value update)
$wsOverview.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Range("G2").Value = "2016-08-28 08:57:20" # Latest HO Xliff Generate Date

$wsDeDe.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H2").Value     = "2016-08-28 08:57:20" # Latest Handoff Datetime (de-de)

$wsZhCn.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H2").Value     = "2016-08-28 08:57:15" # Latest Handoff Datetime (zh-cn)

# --- Widen the columns that now hold "Ready for handoff" ---------------
# Re-fit width (closest achievable via ColumnWidth, which Excel stores in
# 1/6-character increments) to roughly 17.22 characters of stored width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.33  # column E
$wsOverview.Columns.Item(6).ColumnWidth = 16.33  # column F
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33      # column C
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33      # column C
